$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B holds purely numeric strings (phone-style account numbers).
# Format the target range as Text first so Excel stores them as strings
# instead of silently coercing them to numeric values.
$ws.Range("B245:B305").NumberFormat = "@"

$ws.Cells.Item(245, 1).Value = '2026-02-11 14:15:16'
$ws.Cells.Item(245, 2).Value = '237679087694'
$ws.Cells.Item(245, 3).Value = 'RODELPHA TANE TATSIDA'
$ws.Cells.Item(245, 4).Value = 82791

$ws.Cells.Item(246, 1).Value = '2026-02-11 15:44:42'
$ws.Cells.Item(246, 2).Value = '237679869809'
$ws.Cells.Item(246, 3).Value = 'XAVIEE ROSINE MEZAMO'
$ws.Cells.Item(246, 4).Value = 16393

$ws.Cells.Item(247, 1).Value = '2026-02-11 16:31:26'
$ws.Cells.Item(247, 2).Value = '237680435802'
$ws.Cells.Item(247, 3).Value = 'derice mboumela'
$ws.Cells.Item(247, 4).Value = 1525642

$ws.Cells.Item(248, 1).Value = '2026-02-11 13:43:05'
$ws.Cells.Item(248, 2).Value = '237681114247'
$ws.Cells.Item(248, 3).Value = 'LIVIE CHRISTIANE NGOUFACK SONTIA'
$ws.Cells.Item(248, 4).Value = 85354

$ws.Cells.Item(249, 1).Value = '2026-02-11 09:01:34'
$ws.Cells.Item(249, 2).Value = '237681602244'
$ws.Cells.Item(249, 3).Value = 'TSOMEJIO KENFACK NICAISE NESLIE ETS TCHATCHOUANG PAUL  ETP'
$ws.Cells.Item(249, 4).Value = 10339

$ws.Cells.Item(250, 1).Value = '2026-02-11 15:02:37'
$ws.Cells.Item(250, 2).Value = '237681606646'
$ws.Cells.Item(250, 3).Value = 'LA NEGRESSE SARL DONFACK PAULINE PELAGIE'
$ws.Cells.Item(250, 4).Value = 152371

$ws.Cells.Item(251, 1).Value = '2026-02-11 15:28:54'
$ws.Cells.Item(251, 2).Value = '237681655241'
$ws.Cells.Item(251, 3).Value = 'LA NEGRESSE LTDLA CBOX R1 NKUIDJEU KAMDOUM SYMPHORIEN'
$ws.Cells.Item(251, 4).Value = 38042

$ws.Cells.Item(252, 1).Value = '2026-02-11 15:58:37'
$ws.Cells.Item(252, 2).Value = '237681658403'
$ws.Cells.Item(252, 3).Value = 'LA NEGRESSE MISSOKE-UNIVERSITE'
$ws.Cells.Item(252, 4).Value = 147230

$ws.Cells.Item(253, 1).Value = '2026-02-11 14:21:30'
$ws.Cells.Item(253, 2).Value = '237681862876'
$ws.Cells.Item(253, 3).Value = 'TIDO GARLINE NOGRA-POLAS-BTQ-MAKEPE MISSOKE'
$ws.Cells.Item(253, 4).Value = 12408

$ws.Cells.Item(254, 1).Value = '2026-02-11 12:31:35'
$ws.Cells.Item(254, 2).Value = '237683555873'
$ws.Cells.Item(254, 3).Value = 'CHI MERCY SWIRI LTDLA_POLAS_BTQ_LIMBE'
$ws.Cells.Item(254, 4).Value = 658042

$ws.Cells.Item(255, 1).Value = '2026-02-11 15:20:15'
$ws.Cells.Item(255, 2).Value = '237654037914'
$ws.Cells.Item(255, 3).Value = 'LA NEGRESSE SARL YONGA RUSSEL DONALD'
$ws.Cells.Item(255, 4).Value = 24303

$ws.Cells.Item(256, 1).Value = '2026-02-10 10:49:58'
$ws.Cells.Item(256, 2).Value = '237674243367'
$ws.Cells.Item(256, 3).Value = 'EMEGNI NGUEKAM DESTO WILFRIED CHIC MOBILE SARL'
$ws.Cells.Item(256, 4).Value = 229

$ws.Cells.Item(257, 1).Value = '2026-02-11 16:37:29'
$ws.Cells.Item(257, 2).Value = '237681656314'
$ws.Cells.Item(257, 3).Value = 'SWIRRI AZINWI NGANG LA NEGRESSE SARL'
$ws.Cells.Item(257, 4).Value = 448279

$ws.Cells.Item(258, 1).Value = '2026-02-11 05:00:54'
$ws.Cells.Item(258, 2).Value = '237682511457'
$ws.Cells.Item(258, 3).Value = 'FRANFORETTE NWOGUEP KODJOUO'
$ws.Cells.Item(258, 4).Value = 80586

$ws.Cells.Item(259, 1).Value = '2026-02-11 16:36:39'
$ws.Cells.Item(259, 2).Value = '237650934256'
$ws.Cells.Item(259, 3).Value = 'DIDIER ROMUALD MBAKOP NYA'
$ws.Cells.Item(259, 4).Value = 476

$ws.Cells.Item(260, 1).Value = '2026-02-11 14:52:52'
$ws.Cells.Item(260, 2).Value = '237652275301'
$ws.Cells.Item(260, 3).Value = 'NDAMI EPSE NONGA ROSALIE ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Cells.Item(260, 4).Value = 2792

$ws.Cells.Item(261, 1).Value = '2026-02-11 12:29:56'
$ws.Cells.Item(261, 2).Value = '237652427111'
$ws.Cells.Item(261, 3).Value = 'CELESTINE CHANTAL MENDJOBOU EPSE NZIAKOU NJANJO'
$ws.Cells.Item(261, 4).Value = 46749

$ws.Cells.Item(262, 1).Value = '2026-02-11 08:32:49'
$ws.Cells.Item(262, 2).Value = '237671645947'
$ws.Cells.Item(262, 3).Value = 'CLOVIS TAMKOU SOCTOUO'
$ws.Cells.Item(262, 4).Value = 24776

$ws.Cells.Item(263, 1).Value = '2026-02-11 16:03:09'
$ws.Cells.Item(263, 2).Value = '237672064755'
$ws.Cells.Item(263, 3).Value = 'KENFACK FRANC DUVIAL LA NEGRESSE SARL'
$ws.Cells.Item(263, 4).Value = 1155719

$ws.Cells.Item(264, 1).Value = '2026-02-11 11:30:18'
$ws.Cells.Item(264, 2).Value = '237672956746'
$ws.Cells.Item(264, 3).Value = 'ALAIN MOISE NDJONG ITALEN'
$ws.Cells.Item(264, 4).Value = 140

$ws.Cells.Item(265, 1).Value = '2026-02-11 15:14:10'
$ws.Cells.Item(265, 2).Value = '237673718583'
$ws.Cells.Item(265, 3).Value = 'Soppi Verole ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Cells.Item(265, 4).Value = 250482

$ws.Cells.Item(266, 1).Value = '2026-02-11 09:14:48'
$ws.Cells.Item(266, 2).Value = '237674933048'
$ws.Cells.Item(266, 3).Value = 'LA NEGRESSE LTDLA CBOX R1 NINZEGA EMILIENNE CLAIRE'
$ws.Cells.Item(266, 4).Value = 868

$ws.Cells.Item(267, 1).Value = '2026-02-11 14:13:10'
$ws.Cells.Item(267, 2).Value = '237675950748'
$ws.Cells.Item(267, 3).Value = 'ATANGANA HENRI EITEL ETS TCHATCHOUANG PAUL _ETP'
$ws.Cells.Item(267, 4).Value = 222843

$ws.Cells.Item(268, 1).Value = '2026-02-11 12:48:59'
$ws.Cells.Item(268, 2).Value = '237682520113'
$ws.Cells.Item(268, 3).Value = 'KEUYAP NGATCHEU JUDITH JOSY ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Cells.Item(268, 4).Value = 205235

$ws.Cells.Item(269, 1).Value = '2026-02-11 15:27:24'
$ws.Cells.Item(269, 2).Value = '237683356603'
$ws.Cells.Item(269, 3).Value = 'MFS MATCHINDA SENDRINE'
$ws.Cells.Item(269, 4).Value = 11827

$ws.Cells.Item(270, 1).Value = '2026-02-11 16:48:42'
$ws.Cells.Item(270, 2).Value = '237683394976'
$ws.Cells.Item(270, 3).Value = 'DIALLO MAMADOU OURY'
$ws.Cells.Item(270, 4).Value = 13786

$ws.Cells.Item(271, 1).Value = '2026-02-11 15:27:53'
$ws.Cells.Item(271, 2).Value = '237683395123'
$ws.Cells.Item(271, 3).Value = 'LA NEGRESSE SARL LTDLA-CBOX-R1-TAGNIN NICAISSE FLEURIE'
$ws.Cells.Item(271, 4).Value = 105285

$ws.Cells.Item(272, 1).Value = '2026-02-11 16:55:38'
$ws.Cells.Item(272, 2).Value = '237670174030'
$ws.Cells.Item(272, 3).Value = 'LA NEGRESSE LTDLA_CBOX_R1_FOHOM STEPHANE THIERRY'
$ws.Cells.Item(272, 4).Value = 86097

$ws.Cells.Item(273, 1).Value = '2026-02-11 16:28:19'
$ws.Cells.Item(273, 2).Value = '237679111075'
$ws.Cells.Item(273, 3).Value = 'TCHAMABE YOSSA JEAN JOEL ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Cells.Item(273, 4).Value = 289231

$ws.Cells.Item(274, 1).Value = '2026-02-11 17:38:42'
$ws.Cells.Item(274, 2).Value = '237681611433'
$ws.Cells.Item(274, 3).Value = 'Barry Diakariaou World T Plus'
$ws.Cells.Item(274, 4).Value = 73806

$ws.Cells.Item(275, 1).Value = '2026-02-11 12:52:47'
$ws.Cells.Item(275, 2).Value = '237652194260'
$ws.Cells.Item(275, 3).Value = 'NADEGE MALEUTCHOUA'
$ws.Cells.Item(275, 4).Value = 739368

$ws.Cells.Item(276, 1).Value = '2026-02-11 14:12:57'
$ws.Cells.Item(276, 2).Value = '237652667691'
$ws.Cells.Item(276, 3).Value = 'OLIVIA BI'
$ws.Cells.Item(276, 4).Value = 36045

$ws.Cells.Item(277, 1).Value = '2026-02-11 13:39:37'
$ws.Cells.Item(277, 2).Value = '237670473852'
$ws.Cells.Item(277, 3).Value = 'LA NEGRESSE LTDLA CBOX R0 MAGNE TALLA EMILIE'
$ws.Cells.Item(277, 4).Value = 387450

$ws.Cells.Item(278, 1).Value = '2026-02-11 16:13:02'
$ws.Cells.Item(278, 2).Value = '237673018936'
$ws.Cells.Item(278, 3).Value = 'EJUH AKEP EUGENE ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Cells.Item(278, 4).Value = 51934

$ws.Cells.Item(279, 1).Value = '2026-02-11 14:00:26'
$ws.Cells.Item(279, 2).Value = '237673560726'
$ws.Cells.Item(279, 3).Value = 'MAGUELON NADERGE -CHIC MOBILE'
$ws.Cells.Item(279, 4).Value = 9875

$ws.Cells.Item(280, 1).Value = '2026-02-11 16:03:33'
$ws.Cells.Item(280, 2).Value = '237674450580'
$ws.Cells.Item(280, 3).Value = 'TCHIYADJE VANESSA ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Cells.Item(280, 4).Value = 31259

$ws.Cells.Item(281, 1).Value = '2026-02-11 15:48:24'
$ws.Cells.Item(281, 2).Value = '237674929417'
$ws.Cells.Item(281, 3).Value = 'HORTANCE MANTHO'
$ws.Cells.Item(281, 4).Value = 54522

$ws.Cells.Item(282, 1).Value = '2026-02-11 16:01:57'
$ws.Cells.Item(282, 2).Value = '237677316351'
$ws.Cells.Item(282, 3).Value = 'QUELIE LOVE KAGE LEUYOUM'
$ws.Cells.Item(282, 4).Value = 3541

$ws.Cells.Item(283, 1).Value = '2026-02-11 13:17:07'
$ws.Cells.Item(283, 2).Value = '237677831340'
$ws.Cells.Item(283, 3).Value = 'YEPCHE NGANSOP DORINELLE FLORE ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Cells.Item(283, 4).Value = 27133

$ws.Cells.Item(284, 1).Value = '2026-02-11 14:33:54'
$ws.Cells.Item(284, 2).Value = '237678843959'
$ws.Cells.Item(284, 3).Value = 'ETS MOBILE FINANCIAL SERVICES MFS NGO NDJAYICK THERESE'
$ws.Cells.Item(284, 4).Value = 35994

$ws.Cells.Item(285, 1).Value = '2026-02-11 12:22:44'
$ws.Cells.Item(285, 2).Value = '237679068456'
$ws.Cells.Item(285, 3).Value = 'NDATSE EPSE NANWO ANGELE SOLANGE ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Cells.Item(285, 4).Value = 228147

$ws.Cells.Item(286, 1).Value = '2026-02-11 14:57:10'
$ws.Cells.Item(286, 2).Value = '237679793647'
$ws.Cells.Item(286, 3).Value = 'ETS LE CONTENT 71'
$ws.Cells.Item(286, 4).Value = 95

$ws.Cells.Item(287, 1).Value = '2026-02-11 09:48:27'
$ws.Cells.Item(287, 2).Value = '237682323406'
$ws.Cells.Item(287, 3).Value = 'ETS LE CONTENT LAS VEGAS'
$ws.Cells.Item(287, 4).Value = 36

$ws.Cells.Item(288, 1).Value = '2026-02-11 14:19:55'
$ws.Cells.Item(288, 2).Value = '237682764368'
$ws.Cells.Item(288, 3).Value = 'LA NEGRESSE SARL LIEDJI GINETTE'
$ws.Cells.Item(288, 4).Value = 755746

$ws.Cells.Item(289, 1).Value = '2026-02-11 15:46:01'
$ws.Cells.Item(289, 2).Value = '237682814055'
$ws.Cells.Item(289, 3).Value = 'SORELLE REINE MAKANKEU TENE'
$ws.Cells.Item(289, 4).Value = 31910

$ws.Cells.Item(290, 1).Value = '2026-02-11 14:57:32'
$ws.Cells.Item(290, 2).Value = '237683079541'
$ws.Cells.Item(290, 3).Value = 'THIERRY MELINGUI AYINA'
$ws.Cells.Item(290, 4).Value = 243401

$ws.Cells.Item(291, 1).Value = '2026-02-11 14:17:27'
$ws.Cells.Item(291, 2).Value = '237683379070'
$ws.Cells.Item(291, 3).Value = 'MELI DOUANLA ORNELA LINDA ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Cells.Item(291, 4).Value = 806983

$ws.Cells.Item(292, 1).Value = '2026-02-11 13:19:01'
$ws.Cells.Item(292, 2).Value = '237683730580'
$ws.Cells.Item(292, 3).Value = 'DANIELLA KOMGUEP KOUAMO'
$ws.Cells.Item(292, 4).Value = 46349

$ws.Cells.Item(293, 1).Value = '2026-02-11 14:56:33'
$ws.Cells.Item(293, 2).Value = '237651213730'
$ws.Cells.Item(293, 3).Value = 'NZONDE GABRIEL ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Cells.Item(293, 4).Value = 147892

$ws.Cells.Item(294, 1).Value = '2026-02-11 17:12:58'
$ws.Cells.Item(294, 2).Value = '237652071114'
$ws.Cells.Item(294, 3).Value = 'LA NEGRESSE SARL TALLA BOYOM JEAN ERNEST'
$ws.Cells.Item(294, 4).Value = 145243

$ws.Cells.Item(295, 1).Value = '2026-02-11 13:37:54'
$ws.Cells.Item(295, 2).Value = '237652285489'
$ws.Cells.Item(295, 3).Value = 'joseline kenne'
$ws.Cells.Item(295, 4).Value = 108119

$ws.Cells.Item(296, 1).Value = '2026-02-11 17:55:52'
$ws.Cells.Item(296, 2).Value = '237652940152'
$ws.Cells.Item(296, 3).Value = 'GISAWO AIME LE CONTENT'
$ws.Cells.Item(296, 4).Value = 187901

$ws.Cells.Item(297, 1).Value = '2026-02-11 10:30:14'
$ws.Cells.Item(297, 2).Value = '237654137136'
$ws.Cells.Item(297, 3).Value = 'NGANGUE NDOUMBE CHARLOTTE MINDEM SARL'
$ws.Cells.Item(297, 4).Value = 259628

$ws.Cells.Item(298, 1).Value = '2026-02-11 09:58:52'
$ws.Cells.Item(298, 2).Value = '237671694408'
$ws.Cells.Item(298, 3).Value = 'LOUISE STEPHANIE ZEH'
$ws.Cells.Item(298, 4).Value = 144302

$ws.Cells.Item(299, 1).Value = '2026-02-11 10:51:12'
$ws.Cells.Item(299, 2).Value = '237672279571'
$ws.Cells.Item(299, 3).Value = 'MOTCHUENG MADO LIONNELLE_ TOP MOBIL'
$ws.Cells.Item(299, 4).Value = 1477567

$ws.Cells.Item(300, 1).Value = '2026-02-11 12:48:35'
$ws.Cells.Item(300, 2).Value = '237672777139'
$ws.Cells.Item(300, 3).Value = 'CHOURUPOUO MBAKOP ABDEL MFS'
$ws.Cells.Item(300, 4).Value = 165980

$ws.Cells.Item(301, 1).Value = '2026-02-11 13:40:28'
$ws.Cells.Item(301, 2).Value = '237673220938'
$ws.Cells.Item(301, 3).Value = 'KAMGANG TOMDJIO SIMEONE BABETTE ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Cells.Item(301, 4).Value = 3824

$ws.Cells.Item(302, 1).Value = '2026-02-11 17:01:17'
$ws.Cells.Item(302, 2).Value = '237674437082'
$ws.Cells.Item(302, 3).Value = 'MODESTE NGOUKOUA'
$ws.Cells.Item(302, 4).Value = 55319

$ws.Cells.Item(303, 1).Value = '2026-02-11 14:53:52'
$ws.Cells.Item(303, 2).Value = '237674895877'
$ws.Cells.Item(303, 3).Value = 'jeanne annie ngo mback'
$ws.Cells.Item(303, 4).Value = 365789

$ws.Cells.Item(304, 1).Value = '2026-02-11 15:11:48'
$ws.Cells.Item(304, 2).Value = '237674956331'
$ws.Cells.Item(304, 3).Value = 'ODETTE KUYUKEH'
$ws.Cells.Item(304, 4).Value = 256053

$ws.Cells.Item(305, 1).Value = '2026-02-11 10:34:41'
$ws.Cells.Item(305, 2).Value = '237674979451'
$ws.Cells.Item(305, 3).Value = 'EMPIRE COMPANY   LIMITED SWKBA CBOX R4 TCHOUALA GLADIS NADEGE'
$ws.Cells.Item(305, 4).Value = 72051
